$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 11")

# --- Update cell values (times/ids refreshed by the author) ---
$ws.Range("C9").Value  = "3045987650"
$ws.Range("D9").Value  = "732111193278858"

$ws.Range("B10").Value = "484303795"
$ws.Range("C10").Value = "3046010569"
$ws.Range("D10").Value = "732111193280551"
$ws.Range("E10").Value = "3046008593"

$ws.Range("C11").Value = "3045987650"
$ws.Range("D11").Value = "732111193278858"

$ws.Range("C12").Value = "3052749177"
$ws.Range("D12").Value = "732111324709512"

$ws.Range("C13").Value = "3046010523"
$ws.Range("D13").Value = "732111193280544"

$ws.Range("C14").Value = "3045984556"

# --- New row with an (empty) selected cell ---
$ws.Range("B17").Value = ""

# --- Selection / view housekeeping to mirror author's final state ---
$ws.Range("B17").Select()
